# Append-refresh: 2025-10-21 06:28 JST
# Replaces rows 2-9 of the "ランサーズ" sheet with freshly scraped cases,
# drops old rows 10-26, fixes up hyperlinks and shrinks a few column widths.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1) Drop the old trailing rows (10-26) so only the header + 8 data rows
#    remain (dimension becomes A1:H9).
# ---------------------------------------------------------------------
$ws.Rows("10:26").Delete()

# ---------------------------------------------------------------------
# 2) Overwrite rows 2-9 with the newly scraped case information.
# ---------------------------------------------------------------------
$timestamp = "2025-10-21 06:28:04"
$category = "システム開発"
$deadline = "期限情報なし"

$rows = @(
    @{ Row = 2;  Title = "【急募】ebayAPIを活用したShippingポリシー設定の専門家募集"; Price = "20,000 円 ~ 50,000 円 / 固定"; Url = "https://www.lancers.jp/work/detail/5415908"; Score = 183; Skill = "🔥API" },
    @{ Row = 3;  Title = "【 急募! 】 JS、PHPを使用したWEBシステムの開発、修正の対応"; Price = "500,000 円 ~ 1,000,000 円 / 固定"; Url = "https://www.lancers.jp/work/detail/5417295"; Score = 110; Skill = "◆開発 ○PHP" },
    @{ Row = 4;  Title = "システムの開発補助や運営サポート【フルリモート×長期】"; Price = "300,000 円 ~ 500,000 円 / 固定"; Url = "https://www.lancers.jp/work/detail/5408664"; Score = 90;  Skill = "◆開発" },
    @{ Row = 5;  Title = "MySQLバージョンアップ(ロリポップ/WordPress/1データベースに8サイト)"; Price = "20,000 円 ~ 50,000 円 / 固定"; Url = "https://www.lancers.jp/work/detail/5417433"; Score = 88;  Skill = "◇MySQL ○WordPress" },
    @{ Row = 6;  Title = "【急募】MT4/MT5用FX自動売買システムの開発者募集"; Price = "200,000 円 ~ 300,000 円 / 固定"; Url = "https://www.lancers.jp/work/detail/5417377"; Score = 83;  Skill = "◆開発" },
    @{ Row = 7;  Title = "【急募】Salesforce・MA・CRMコンサルタント経験者を探しています!"; Price = "200,000 円 ~ 300,000 円 / 固定"; Url = "https://www.lancers.jp/work/detail/5371747"; Score = 48;  Skill = "◆コンサル" },
    @{ Row = 8;  Title = "【高齢者支援】見守りアプリとマニュアル制作の依頼"; Price = "300,000 円 ~ 500,000 円 / 固定"; Url = "https://www.lancers.jp/work/detail/5417267"; Score = 45;  Skill = "◇アプリ" },
    @{ Row = 9;  Title = "【 急募】コミュニティサイトのカスタマイズ、修正"; Price = "50,000 円 ~ 100,000 円 / 固定"; Url = "https://www.lancers.jp/work/detail/5417308"; Score = 38;  Skill = "◇サイト" }
)

foreach ($item in $rows) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = $timestamp
    $ws.Cells.Item($r, 2).Value = $item.Title
    $ws.Cells.Item($r, 3).Value = $category
    $ws.Cells.Item($r, 4).Value = $item.Price
    $ws.Cells.Item($r, 5).Value = $deadline
    $ws.Cells.Item($r, 6).Value = $item.Url
    $ws.Cells.Item($r, 7).Value = $item.Score
    $ws.Cells.Item($r, 8).Value = $item.Skill
}

# ---------------------------------------------------------------------
# 3) Rebuild the hyperlinks so they point at the refreshed URLs.
#    (Deleting any single hyperlink clears the whole collection in this
#    engine, so do it once up-front and then re-add the ones we need.)
# ---------------------------------------------------------------------
$ws.Range("F2").Hyperlinks.Delete()

foreach ($item in $rows) {
    $r = $item.Row
    $ws.Hyperlinks.Add($ws.Cells.Item($r, 6), $item.Url) | Out-Null
}

# ---------------------------------------------------------------------
# 4) Shrink a few column widths (stored width = ColumnWidth + 5/6).
# ---------------------------------------------------------------------
$padding = 5 / 6
$ws.Columns.Item(2).ColumnWidth = 45 - $padding
$ws.Columns.Item(4).ColumnWidth = 30 - $padding
$ws.Columns.Item(8).ColumnWidth = 19 - $padding
